$d = $word.ActiveDocument

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$rA = $lastPara.Range.Duplicate
$rA.Collapse(1)
$xmlA = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titolo1"/></w:pPr><w:r><w:t xml:space="preserve">Versione senza </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>contact</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Titolo2"/></w:pPr><w:r><w:t>Sensore</w:t></w:r></w:p><w:p><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">Nel metodo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>sendData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, anziché mandare i dati tramite </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>contact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> fai sì che ogni mezzo secondo venga invocato il metodo  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>riceviDatiSensore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>SC</w:t></w:r><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, anziché mandare i dati tramite </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>contact</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Titolo2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>UserCmd</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">nel metodo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>MandaComando</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Scontrol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> anziché il metodo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>userCmdDemand</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> si deve invocare il metodo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>riceviEdElaboraComandoUserCmdFT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Scontrol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> e memorizzare lo status restituito</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titolo2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Scontrol</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">Puoi commentare ciò che avveniva nei metodi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Dojob</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>run</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>thread</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, visto </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>cheScontrol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> non deve più attendere </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>mesaggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>contact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> provenienti da Sensori e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>UserCmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>, ma riceve solo normali chiamate</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titolo2"/></w:pPr><w:r><w:t>Riferimenti e configurazione</w:t></w:r></w:p><w:p><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">Sia i sensori che </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>UserCmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> devono avere </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>un’attributo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Scontrol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, con relativo setter. Il riferimento dei sensori viene settato nel metodo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>configure</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">  della classe Edi, mentre quello di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>userCmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> viene fatto all’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>internoo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> del metodo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>getInstance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> che ha come </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>paramentro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>userCmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Scontrol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>. Lasci comunque invariato in Edi l’assegnazione di nomi ai processi e il loro avvio</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rA.InsertXML($xmlA)

$lastPara2 = $d.Paragraphs($d.Paragraphs.Count)
$rB = $lastPara2.Range.Duplicate
$rB.Collapse(1)
$xmlB = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titolo1"/></w:pPr><w:r><w:t>Interfaccia utente</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>UserCmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">  deve, tramite il pattern </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>observer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, notificare l’interfaccia utente quando il suo status cambia. Quindi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>UserCmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> deve implementare i metodi ( alla </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>observer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>addGui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>NorifyGui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Remove</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Gui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">. I primi due aggiungono o tolgono ad un vettore di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>gui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> un oggetto che implementa  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>IGui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> il terzo chiama i metodi update delle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>gui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve"> passando come parametro lo status .</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rB.InsertXML($xmlB)

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
